$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the scraped crypto price/volume(1h) data (and, for a couple of
# rows whose rank order changed, the coin name/link too).
# Column layout: B=Coin, C=Link, D=Price, E=Volume(1h)
# D/E are stored as text (not numbers) in the sheet, so force the "Text"
# number format before writing -- otherwise Excel would silently reinterpret
# values like "0.9990" or "29.933.31" as numbers/dates and mangle them.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.933.31'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.13%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.899.66'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.72%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9992'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.10%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7603'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +2.74%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '240.69'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.18%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9996'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.01%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3056'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -2.83%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '25.46'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -6.36%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06851'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -2.17%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07995'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.37%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7529'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -3.77%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.900.25'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.69%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.231'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.39%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.31'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.84%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '29.927.59'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.07%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.96'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.84%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.972'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +2.26%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '241.36'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.62%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007739'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.63%  '

$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.162.11'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.86%  '

$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9990'

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9990'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.11%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.995'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +4.84%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.260'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -2.42%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '165.59'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.36%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.77'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.28%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1307'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +2.49%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.032'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -4.28%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.376'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +2.03%  '

$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.77%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.301'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.82%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.037'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.21%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05364'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.89%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.257'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -3.90%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7298'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -2.99%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.720'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -1.35%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01929'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.97%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.777'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.83%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.195'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.94%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4425'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.73%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '72.46'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -4.62%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.915'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.58%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9996'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.04%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8296'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.45%  '

$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.604'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.03%  '

$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '101.10'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.13%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.777'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.24%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.061.29'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.90%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '36.23'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -3.14%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05956'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.73%  '
